$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4115256666666666
$ws.Range("H2").Value = 1.234577
$ws.Range("I2").Value = 0.2245998342667577
$ws.Range("J2").Value = 0.2245998342667577
$ws.Range("M2").Value = 110.642708
$ws.Range("N2").Value = 331.928124
$ws.Range("O2").Value = 0.5476418925386564
$ws.Range("P2").Value = 0.5476418925386564
$ws.Range("Q2").Value = 45.53231417150532
$ws.Range("R2").Value = 409.7908275435479
$ws.Range("S2").Value = 0.1230002783017158
$ws.Range("T2").Value = 0.1230002783017158
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4115256666666666
$ws.Range("H3").Value = 1.234577
$ws.Range("I3").Value = 0.2245998342667577
$ws.Range("J3").Value = 0.2245998342667577
$ws.Range("O3").Value = 0.3151072754333865
$ws.Range("P3").Value = 0.3151072754333865
$ws.Range("Q3").Value = 26.19880556662722
$ws.Range("R3").Value = 235.789250099645
$ws.Range("S3").Value = 0.07077304183858819
$ws.Range("T3").Value = 0.07077304183858819
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.4115256666666666
$ws.Range("H4").Value = 1.234577
$ws.Range("I4").Value = 0.2245998342667577
$ws.Range("J4").Value = 0.2245998342667577
$ws.Range("M4").Value = 27.72944133333333
$ws.Range("N4").Value = 83.18832399999999
$ws.Range("O4").Value = 0.1372508320279571
$ws.Range("P4").Value = 0.1372508320279571
$ws.Range("Q4").Value = 11.41137683099422
$ws.Range("R4").Value = 102.702391478948
$ws.Range("S4").Value = 0.03082651412645377
$ws.Range("T4").Value = 0.03082651412645377
$ws.Range("I5").Value = 0.3944722233087159
$ws.Range("J5").Value = 0.3944722233087159
$ws.Range("M5").Value = 110.642708
$ws.Range("N5").Value = 331.928124
$ws.Range("O5").Value = 0.5476418925386564
$ws.Range("P5").Value = 0.5476418925386564
$ws.Range("Q5").Value = 79.96993079831066
$ws.Range("R5").Value = 719.7293771847959
$ws.Range("S5").Value = 0.2160295149267167
$ws.Range("T5").Value = 0.2160295149267167
$ws.Range("I6").Value = 0.3944722233087159
$ws.Range("J6").Value = 0.3944722233087159
$ws.Range("O6").Value = 0.3151072754333865
$ws.Range("P6").Value = 0.3151072754333865
$ws.Range("S6").Value = 0.1243010675209599
$ws.Range("T6").Value = 0.1243010675209599
$ws.Range("I7").Value = 0.3944722233087159
$ws.Range("J7").Value = 0.3944722233087159
$ws.Range("M7").Value = 27.72944133333333
$ws.Range("N7").Value = 83.18832399999999
$ws.Range("O7").Value = 0.1372508320279571
$ws.Range("P7").Value = 0.1372508320279571
$ws.Range("Q7").Value = 20.04218393228844
$ws.Range("R7").Value = 180.379655390596
$ws.Range("S7").Value = 0.05414164086103936
$ws.Range("T7").Value = 0.05414164086103936
$ws.Range("G8").Value = 0.6979596666666668
$ws.Range("H8").Value = 2.093879
$ws.Range("I8").Value = 0.3809279424245264
$ws.Range("J8").Value = 0.3809279424245264
$ws.Range("M8").Value = 110.642708
$ws.Range("N8").Value = 331.928124
$ws.Range("O8").Value = 0.5476418925386564
$ws.Range("P8").Value = 0.5476418925386564
$ws.Range("Q8").Value = 77.22414759477734
$ws.Range("R8").Value = 695.0173283529961
$ws.Range("S8").Value = 0.208612099310224
$ws.Range("T8").Value = 0.208612099310224
$ws.Range("G9").Value = 0.6979596666666668
$ws.Range("H9").Value = 2.093879
$ws.Range("I9").Value = 0.3809279424245264
$ws.Range("J9").Value = 0.3809279424245264
$ws.Range("O9").Value = 0.3151072754333865
$ws.Range("P9").Value = 0.3151072754333865
$ws.Range("Q9").Value = 44.43394685065723
$ws.Range("R9").Value = 399.905521655915
$ws.Range("S9").Value = 0.1200331660738384
$ws.Range("T9").Value = 0.1200331660738384
$ws.Range("G10").Value = 0.6979596666666668
$ws.Range("H10").Value = 2.093879
$ws.Range("I10").Value = 0.3809279424245264
$ws.Range("J10").Value = 0.3809279424245264
$ws.Range("M10").Value = 27.72944133333333
$ws.Range("N10").Value = 83.18832399999999
$ws.Range("O10").Value = 0.1372508320279571
$ws.Range("P10").Value = 0.1372508320279571
$ws.Range("Q10").Value = 19.35403162986622
$ws.Range("R10").Value = 174.186284668796
$ws.Range("S10").Value = 0.05228267704046399
$ws.Range("T10").Value = 0.05228267704046399
